# Update cryptocurrency price/volume figures per the Feb 4 2023 13:00 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.84"
$ws.Range("E2").Value = "'-0.06%"
$ws.Range("D3").Value = "'41.57"
$ws.Range("E3").Value = "'6.27%"
$ws.Range("D4").Value = "'5.688"
$ws.Range("E4").Value = "'-1.02%"
$ws.Range("D5").Value = "'0.08349"
$ws.Range("E5").Value = "'4.06%"
$ws.Range("D6").Value = "'2.037"
$ws.Range("E6").Value = "'5.73%"
$ws.Range("D7").Value = "'8.809"
$ws.Range("E7").Value = "'2.20%"
$ws.Range("D8").Value = "'4.533"
$ws.Range("E8").Value = "'0.52%"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("D10").Value = "'0.9298"
$ws.Range("E10").Value = "'1.33%"
$ws.Range("D11").Value = "'0.1297"
$ws.Range("E11").Value = "'3.15%"
$ws.Range("D12").Value = "'0.1964"
$ws.Range("E12").Value = "'0.97%"
$ws.Range("D13").Value = "'0.09380"
$ws.Range("E13").Value = "'1.62%"
$ws.Range("D14").Value = "'0.03917"
$ws.Range("E14").Value = "'9.60%"
$ws.Range("E15").Value = "'0.77%"
$ws.Range("D16").Value = "'0.001300"
$ws.Range("E16").Value = "'0.06%"
$ws.Range("D17").Value = "'0.006192"
$ws.Range("E17").Value = "'-2.76%"
$ws.Range("E18").Value = "'2.34%"
$ws.Range("E19").Value = "'2.13%"
$ws.Range("D20").Value = "'8.312"
$ws.Range("E20").Value = "'-4.45%"
$ws.Range("D21").Value = "'0.1362"
$ws.Range("E21").Value = "'-0.95%"
$ws.Range("D22").Value = "'0.2406"
$ws.Range("E22").Value = "'-10.12%"
$ws.Range("D23").Value = "'0.04417"
$ws.Range("E23").Value = "'-0.78%"
$ws.Range("E24").Value = "'-1.18%"
$ws.Range("D25").Value = "'0.004380"
$ws.Range("E25").Value = "'-0.99%"
$ws.Range("E26").Value = "'-0.54%"
$ws.Range("D39").Value = "'0.02816"
$ws.Range("E39").Value = "'12.09%"
$ws.Range("D40").Value = "'0.05554"
$ws.Range("E40").Value = "'2.16%"
$ws.Range("D41").Value = "'0.007806"
$ws.Range("E41").Value = "'4.12%"
$ws.Range("E42").Value = "'2.62%"
$ws.Range("D43").Value = "'0.008936"
$ws.Range("E43").Value = "'-10.18%"
$ws.Range("E44").Value = "'1.00%"
$ws.Range("D45").Value = "'0.01175"
$ws.Range("E45").Value = "'3.30%"
$ws.Range("D46").Value = "'0.00007011"
$ws.Range("E46").Value = "'2.63%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.53%"
$ws.Range("D48").Value = "'0.003178"
$ws.Range("E48").Value = "'4.14%"
$ws.Range("D49").Value = "'0.002276"
$ws.Range("E49").Value = "'-0.60%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.53%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.53%"
